$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update volume/date header text (rich-text shared strings -> set plain text value)
$ws.Range("A8").Value = "Volume 32   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# --- Crime statistics table updates ---
$ws.Range("G14").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4163) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4163) | Out-Null
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 20
$ws.Range("H16").Value = 53.846153846153
$ws.Range("I16").Value = 77
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = -11.494252873563
$ws.Range("L16").Value = -10.465116279069
$ws.Range("M16").Value = -9.411764705882
$ws.Range("N16").Value = -85.077519379845
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 68.421052631578
$ws.Range("I17").Value = 162
$ws.Range("J17").Value = 143
$ws.Range("K17").Value = 13.286713286713
$ws.Range("L17").Value = 29.6
$ws.Range("M17").Value = 179.310344827586
$ws.Range("N17").Value = -10.49723756906
$ws.Range("I15").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -46.153846153846
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 89
$ws.Range("K18").Value = -1.123595505617
$ws.Range("L18").Value = -21.428571428571
$ws.Range("M18").Value = -35.294117647058
$ws.Range("N18").Value = -92.407247627264
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 41
$ws.Range("H19").Value = 7.894736842105
$ws.Range("I19").Value = 291
$ws.Range("J19").Value = 349
$ws.Range("K19").Value = -16.618911174785
$ws.Range("L19").Value = -18.941504178273
$ws.Range("M19").Value = 43.349753694581
$ws.Range("N19").Value = -62.692307692307
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -81.818181818181
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -43.478260869565
$ws.Range("I20").Value = 129
$ws.Range("J20").Value = 134
$ws.Range("K20").Value = -3.731343283582
$ws.Range("L20").Value = 21.698113207547
$ws.Range("M20").Value = 20.560747663551
$ws.Range("N20").Value = -93.377823408624
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -31.578947368421
$ws.Range("F21").Value = 114
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = 6.542056074766
$ws.Range("I21").Value = 758
$ws.Range("J21").Value = 812
$ws.Range("K21").Value = -6.650246305418
$ws.Range("L21").Value = -5.25
$ws.Range("M21").Value = 27.609427609427
$ws.Range("N21").Value = -83.510985425277
$ws.Range("I15").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 3
$ws.Range("G14").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4163) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4163) | Out-Null
$ws.Range("I22").Value = 8
$ws.Range("K22").Value = -38.461538461538
$ws.Range("L22").Value = 14.285714285714
$ws.Range("M22").Value = 33.333333333333
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 28.571428571428
$ws.Range("I23").Value = 35
$ws.Range("J23").Value = 43
$ws.Range("K23").Value = -18.60465116279
$ws.Range("L23").Value = 2.941176470588
$ws.Range("M23").Value = 169.230769230769
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -20
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = -33.043478260869
$ws.Range("I24").Value = 604
$ws.Range("J24").Value = 629
$ws.Range("K24").Value = -3.974562798092
$ws.Range("L24").Value = -12.463768115942
$ws.Range("M24").Value = 6.338028169014
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 52
$ws.Range("H25").Value = -48.076923076923
$ws.Range("I25").Value = 231
$ws.Range("J25").Value = 252
$ws.Range("K25").Value = -8.333333333333
$ws.Range("L25").Value = -25.483870967741
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("G26").Value = 47
$ws.Range("H26").Value = -12.765957446808
$ws.Range("I26").Value = 282
$ws.Range("J26").Value = 269
$ws.Range("K26").Value = 4.832713754646
$ws.Range("L26").Value = 26.457399103139
$ws.Range("M26").Value = 3.296703296703
$ws.Range("G14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4163) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4163) | Out-Null
$ws.Range("C28").Value = 2
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4163) | Out-Null
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 29
$ws.Range("K28").Value = 45
$ws.Range("L28").Value = 16
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4163) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4163) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4163) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4163) | Out-Null
$ws.Range("H31").Value = -100
$ws.Range("I15").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4122) | Out-Null
$ws.Range("F33").Value = 1
$ws.Range("I33").Value = 4
$ws.Range("K33").Value = -42.857142857142
$ws.Range("L33").Value = 33.333333333333

$ws.Application.CutCopyMode = $false
$wb.Application.Calculate()
